$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.394.18"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.233.75"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'245.08"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").Value = "'0.628"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "'74.02"
$ws.Range("E7").Value = "  -4.11%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "'0.619"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "'42.97"
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").Value = "'0.0967"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").Value = "'7.11"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").Value = "'14.43"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").Value = "'0.850"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "2.220.53"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "42.218.66"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("E18").Value = "  +13.42%  "
$ws.Range("D19").Value = "'6.16"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "'72.07"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "'10.38"
$ws.Range("E21").Value = "  +37.85%  "
$ws.Range("D22").Value = "'231.03"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "'2.16"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("D24").Value = "'11.71"
$ws.Range("E24").Value = "  +3.45%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'3.68"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("D27").Value = "'2.30"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").Value = "'166.62"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").Value = "'20.92"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").Value = "'5.86"
$ws.Range("E31").Value = "  +18.37%  "
$ws.Range("D32").Value = "'0.0805"
$ws.Range("E32").Value = "  -3.82%  "
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "'0.125"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'29.61"
$ws.Range("E35").Value = "  -9.16%  "
$ws.Range("D36").Value = "'4.44"
$ws.Range("E36").Value = "  -1.74%  "
$ws.Range("D37").Value = "'0.0308"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").Value = "'13.16"
$ws.Range("E38").Value = "  -9.29%  "
$ws.Range("D39").Value = "'2.16"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").Value = "'5.63"
$ws.Range("E40").Value = "  -4.39%  "
$ws.Range("D41").Value = "'63.10"
$ws.Range("E41").Value = "  +3.62%  "
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").Value = "'8.82"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("D44").Value = "'105.29"
$ws.Range("E44").Value = "  -6.60%  "
$ws.Range("E45").Value = "  +3.08%  "
$ws.Range("D46").Value = "'0.996"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'1.13"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.37"
$ws.Range("E48").Value = "  +3.91%  "
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "'4.06"
$ws.Range("E51").Value = "  -2.91%  "
